$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Cells.Item(2, 4)
$r.NumberFormat = "@"
$r.Value = '29.787.07'
$r.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +2.13%  '

$r = $ws.Cells.Item(3, 4)
$r.NumberFormat = "@"
$r.Value = '1.859.95'
$r.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +1.66%  '

$r = $ws.Cells.Item(4, 4)
$r.NumberFormat = "@"
$r.Value = '0.9996'
$r.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$r = $ws.Cells.Item(5, 4)
$r.NumberFormat = "@"
$r.Value = '245.05'
$r.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.92%  '

$r = $ws.Cells.Item(6, 4)
$r.NumberFormat = "@"
$r.Value = '0.6426'
$r.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +4.26%  '

$ws.Cells.Item(7, 5).Value = '  -0.04%  '

$ws.Cells.Item(8, 2).Value = 'Dogecoin'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$r = $ws.Cells.Item(8, 4)
$r.NumberFormat = "@"
$r.Value = '0.07542'
$r.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +2.70%  '

$ws.Cells.Item(9, 2).Value = 'Cardano'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$r = $ws.Cells.Item(9, 4)
$r.NumberFormat = "@"
$r.Value = '0.2979'
$r.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +2.40%  '

$ws.Cells.Item(10, 2).Value = 'Solana'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$r = $ws.Cells.Item(10, 4)
$r.NumberFormat = "@"
$r.Value = '24.60'
$r.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +6.18%  '

$ws.Cells.Item(11, 2).Value = 'TRON'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$r = $ws.Cells.Item(11, 4)
$r.NumberFormat = "@"
$r.Value = '0.07687'
$r.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.65%  '

$ws.Cells.Item(12, 2).Value = 'WrappedEther'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$r = $ws.Cells.Item(12, 4)
$r.NumberFormat = "@"
$r.Value = '1.858.75'
$r.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +1.51%  '

$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$r = $ws.Cells.Item(13, 4)
$r.NumberFormat = "@"
$r.Value = '5.051'
$r.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +1.53%  '

$ws.Cells.Item(14, 2).Value = 'Polygon'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$r = $ws.Cells.Item(14, 4)
$r.NumberFormat = "@"
$r.Value = '0.6931'
$r.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +3.39%  '

$ws.Cells.Item(15, 2).Value = 'Litecoin'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$r = $ws.Cells.Item(15, 4)
$r.NumberFormat = "@"
$r.Value = '84.05'
$r.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +2.07%  '

$ws.Cells.Item(16, 2).Value = 'ShibaInu'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$r = $ws.Cells.Item(16, 4)
$r.NumberFormat = "@"
$r.Value = '0.000009865'
$r.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +9.82%  '

$ws.Cells.Item(17, 2).Value = 'Uniswap'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$r = $ws.Cells.Item(17, 4)
$r.NumberFormat = "@"
$r.Value = '6.146'
$r.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +5.26%  '

$ws.Cells.Item(18, 2).Value = 'WrappedBTC'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$r = $ws.Cells.Item(18, 4)
$r.NumberFormat = "@"
$r.Value = '29.805.64'
$r.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +2.24%  '

$ws.Cells.Item(19, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$r = $ws.Cells.Item(19, 4)
$r.NumberFormat = "@"
$r.Value = '2.115.93'
$r.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +1.97%  '

$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$r = $ws.Cells.Item(20, 4)
$r.NumberFormat = "@"
$r.Value = '237.11'
$r.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.44%  '

$ws.Cells.Item(21, 2).Value = 'Avalanche'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$r = $ws.Cells.Item(21, 4)
$r.NumberFormat = "@"
$r.Value = '12.67'
$r.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +1.52%  '

$ws.Cells.Item(22, 2).Value = 'Dai'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$r = $ws.Cells.Item(22, 4)
$r.NumberFormat = "@"
$r.Value = '1.0000'
$r.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.05%  '

$ws.Cells.Item(23, 2).Value = 'Chainlink'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$r = $ws.Cells.Item(23, 4)
$r.NumberFormat = "@"
$r.Value = '7.527'
$r.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +2.35%  '

$ws.Cells.Item(24, 2).Value = 'BinanceUSD'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$r = $ws.Cells.Item(24, 4)
$r.NumberFormat = "@"
$r.Value = '1.001'
$r.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -0.01%  '

$ws.Cells.Item(25, 2).Value = 'Monero'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$r = $ws.Cells.Item(25, 4)
$r.NumberFormat = "@"
$r.Value = '158.86'
$r.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.18%  '

$ws.Cells.Item(26, 2).Value = 'Stellar'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$r = $ws.Cells.Item(26, 4)
$r.NumberFormat = "@"
$r.Value = '0.1425'
$r.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +2.83%  '

$ws.Cells.Item(27, 2).Value = 'Cosmos'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$r = $ws.Cells.Item(27, 4)
$r.NumberFormat = "@"
$r.Value = '8.560'
$r.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.55%  '

$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$r = $ws.Cells.Item(28, 4)
$r.NumberFormat = "@"
$r.Value = '17.93'
$r.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +1.82%  '

$ws.Cells.Item(29, 2).Value = 'Hedera'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$r = $ws.Cells.Item(29, 4)
$r.NumberFormat = "@"
$r.Value = '0.06222'
$r.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +5.93%  '

$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$r = $ws.Cells.Item(30, 4)
$r.NumberFormat = "@"
$r.Value = '1.497'
$r.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.68%  '

$ws.Cells.Item(31, 2).Value = 'Toncoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$r = $ws.Cells.Item(31, 4)
$r.NumberFormat = "@"
$r.Value = '1.291'
$r.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +5.51%  '

$ws.Cells.Item(32, 2).Value = 'Filecoin'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$r = $ws.Cells.Item(32, 4)
$r.NumberFormat = "@"
$r.Value = '4.155'
$r.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +1.78%  '

$ws.Cells.Item(33, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$r = $ws.Cells.Item(33, 4)
$r.NumberFormat = "@"
$r.Value = '4.107'
$r.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +0.72%  '

$ws.Cells.Item(34, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$r = $ws.Cells.Item(34, 4)
$r.NumberFormat = "@"
$r.Value = '1.903'
$r.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +2.52%  '

$ws.Cells.Item(35, 2).Value = 'ARBITRUM'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$r = $ws.Cells.Item(35, 4)
$r.NumberFormat = "@"
$r.Value = '1.174'
$r.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +3.34%  '

$ws.Cells.Item(36, 2).Value = 'ImmutableX'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$r = $ws.Cells.Item(36, 4)
$r.NumberFormat = "@"
$r.Value = '0.7287'
$r.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +0.72%  '

$ws.Cells.Item(37, 2).Value = 'HuobiToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$r = $ws.Cells.Item(37, 4)
$r.NumberFormat = "@"
$r.Value = '2.608'
$r.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -0.27%  '

$ws.Cells.Item(38, 2).Value = 'MXToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$r = $ws.Cells.Item(38, 4)
$r.NumberFormat = "@"
$r.Value = '2.818'
$r.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -1.60%  '

$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$r = $ws.Cells.Item(39, 4)
$r.NumberFormat = "@"
$r.Value = '0.01787'
$r.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +1.64%  '

$ws.Cells.Item(40, 2).Value = 'Maker'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$r = $ws.Cells.Item(40, 4)
$r.NumberFormat = "@"
$r.Value = '1.213.61'
$r.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -1.00%  '

$ws.Cells.Item(41, 2).Value = 'FraxShare'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$r = $ws.Cells.Item(41, 4)
$r.NumberFormat = "@"
$r.Value = '6.320'
$r.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +2.03%  '

$r = $ws.Cells.Item(42, 4)
$r.NumberFormat = "@"
$r.Value = '0.9222'
$r.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +2.07%  '

$ws.Cells.Item(43, 2).Value = 'PaxDollar'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$r = $ws.Cells.Item(43, 4)
$r.NumberFormat = "@"
$r.Value = '1.000'
$r.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.09%  '

$ws.Cells.Item(44, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$r = $ws.Cells.Item(44, 4)
$r.NumberFormat = "@"
$r.Value = '2.029.33'
$r.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.64%  '

$ws.Cells.Item(45, 2).Value = 'Quant'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$r = $ws.Cells.Item(45, 4)
$r.NumberFormat = "@"
$r.Value = '102.13'
$r.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.24%  '

$ws.Cells.Item(46, 2).Value = 'Aave'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$r = $ws.Cells.Item(46, 4)
$r.NumberFormat = "@"
$r.Value = '67.08'
$r.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +2.19%  '

$ws.Cells.Item(47, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$r = $ws.Cells.Item(47, 4)
$r.NumberFormat = "@"
$r.Value = '0.00000000118'
$r.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.66%  '

$ws.Cells.Item(48, 2).Value = 'TheSandbox'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$r = $ws.Cells.Item(48, 4)
$r.NumberFormat = "@"
$r.Value = '0.4060'
$r.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +0.44%  '

$r = $ws.Cells.Item(49, 4)
$r.NumberFormat = "@"
$r.Value = '9.182'
$r.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -0.11%  '

$ws.Cells.Item(50, 2).Value = 'RenderToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$r = $ws.Cells.Item(50, 4)
$r.NumberFormat = "@"
$r.Value = '1.673'
$r.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +5.68%  '

$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$r = $ws.Cells.Item(51, 4)
$r.NumberFormat = "@"
$r.Value = '0.05787'
$r.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.71%  '
